$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add four new rows (11-14) to the LCOE INS table, mirroring rows 7-10 but
# targeting the new "SNK_DAC" sink/DAC technology instead of "EN_Z*".
# ---------------------------------------------------------------------------

# Row 11
$ws.Range("D11").Value = "LO"
$ws.Range("E11").Value = "ACT_BND"
$ws.Range("F11").Value = 2025
$ws.Range("G11").Value = "'-PASTI"
$ws.Range("H11").Borders.LineStyle = -4142
$ws.Range("I11").Value = 0
$ws.Range("J11").Borders.LineStyle = -4142
$ws.Range("K11").Borders.LineStyle = -4142
$ws.Range("L11").Value = "SNK_DAC"

# Row 12
$ws.Range("D12").Value = "LO"
$ws.Range("E12").Value = "ACT_BND"
$ws.Range("F12").Value = 2050
$ws.Range("G12").Value = "'-PASTI"
$ws.Range("H12").Borders.LineStyle = -4142
$ws.Range("I12").Value = 0.005
$ws.Range("J12").Borders.LineStyle = -4142
$ws.Range("K12").Borders.LineStyle = -4142
$ws.Range("L12").Value = "SNK_DAC"

# Row 13
$ws.Range("D13").Value = "LO"
$ws.Range("E13").Value = "ACT_BND"
$ws.Range("F13").Value = 2080
$ws.Range("G13").Value = "'-PASTI"
$ws.Range("H13").Borders.LineStyle = -4142
$ws.Range("I13").Value = 0.01
$ws.Range("J13").Borders.LineStyle = -4142
$ws.Range("K13").Borders.LineStyle = -4142
$ws.Range("L13").Value = "SNK_DAC"

# Row 14
$ws.Range("D14").Value = "LO"
$ws.Range("E14").Value = "ACT_BND"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "'-PASTI"
$ws.Range("H14").Borders.LineStyle = -4142
$ws.Range("I14").Value = 1
$ws.Range("J14").Borders.LineStyle = -4142
$ws.Range("K14").Borders.LineStyle = -4142
$ws.Range("L14").Value = "SNK_DAC"

# ---------------------------------------------------------------------------
# Widen column G so the "-PASTI" values remain fully visible (bestFit-style).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 10

# ---------------------------------------------------------------------------
# Leave a note on the new start-year cell, explaining the 2020 -> 2025 change
# (same note text/author as the existing comment on F7).
# ---------------------------------------------------------------------------
$note = "Mahmoud Mobir:`n12-8-2021`nThis was 2020 but it caused infeasiblities. Made it 2025. "
$ws.Range("F11").AddComment($note) | Out-Null

# ---------------------------------------------------------------------------
# Match the final selection left behind in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("N12").Select() | Out-Null
